# Visium v2.0 - revert to the "old" dataset_type branch:
#   - drop "RNAseq (Visium)", "GeoMx" and "RNAseq (GeoMx)"
#   - add "GeoMx (NGS)" (after CODEX) and "GeoMx (nCounter)" (after 10X Multiome)
#   - dataset_type list shrinks from 35 to 34 entries
#   - bump pav:createdOn on the .metadata sheet

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rewrite the `dataset_type` lookup sheet with the new 34-row list.
# ---------------------------------------------------------------------------
$dsType = $wb.Worksheets.Item("dataset_type")

$labels = @(
    "HiFi-Slide",
    "SNARE-seq2",
    "MIBI",
    "DESI",
    "scATACseq",
    "Auto-fluorescence",
    "Confocal",
    "scRNAseq",
    "Xenium",
    "snATACseq",
    "Molecular Cartography",
    "CosMx",
    "DBiT",
    "SIMS",
    "Cell DIVE",
    "CODEX",
    "GeoMx (NGS)",
    "CyCIF",
    "Light Sheet",
    "RNAseq (bulk)",
    "MALDI",
    "2D Imaging Mass Cytometry",
    "Histology",
    "Enhanced Stimulated Raman Spectroscopy (SRS)",
    "ATACseq (bulk)",
    "MERFISH",
    "LC-MS",
    "10X Multiome",
    "GeoMx (nCounter)",
    "PhenoCycler",
    "Second Harmonic Generation (SHG)",
    "Thick section Multiphoton MxIF",
    "snRNAseq",
    "Visium"
)

$uris = @(
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000195",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000264",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000172",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000204",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000247",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000205",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000206",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000248",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000219",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000183",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000217",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000218",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000222",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000202",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000159",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000160",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000300",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000200",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000168",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000212",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000201",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000296",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000197",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000209",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000210",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000221",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000194",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000215",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000301",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000199",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000208",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000207",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000184",
    "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000187"
)

for ($i = 0; $i -lt $labels.Length; $i++) {
    $row = $i + 1
    $dsType.Cells.Item($row, 1).Value = $labels[$i]
    $dsType.Cells.Item($row, 2).Value = $uris[$i]
}

# The old list had 35 rows; the new one only has 34, so drop the trailing row.
$dsType.Rows.Item(35).Delete()

# ---------------------------------------------------------------------------
# 2. Point the Visium sheet's dataset_type dropdown at the shrunk range.
# ---------------------------------------------------------------------------
$visium = $wb.Worksheets.Item("Visium")
$visium.Range("D2:D1001").Validation.Formula1 = "='dataset_type'!`$A`$1:`$A`$34"

# ---------------------------------------------------------------------------
# 3. Bump the recorded pav:createdOn timestamp on the .metadata sheet.
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item(".metadata")
$meta.Cells.Item(2, 3).Value = "2023-11-15T17:24:29-08:00"
